$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (styles for columns A and B, borders, number format) from the
# last currently-formatted row (30) down through the new last row (62) by copying
# the existing row 30 formatting into the new rows before writing values.
$ws.Range("A30:C30").Copy()
$ws.Range("A31:C62").PasteSpecial(-4122)

# Data for rows 2-62: row number, date serial (column B), hours value (column C).
# Column A always holds (row number - 2).
$data = @(
    @(2, 44347, 1.25466810966811),
    @(3, 44348, 1.194706818648593),
    @(4, 44349, 1.206640023682652),
    @(5, 44350, 0.9256727430555556),
    @(6, 44351, 1.347636973180077),
    @(7, 44352, 1.162026270456503),
    @(8, 44353, 1.044696296296296),
    @(9, 44354, 1.152830394387562),
    @(10, 44355, 1.218148306697108),
    @(11, 44356, 1.222159953298307),
    @(12, 44357, 1.168707034372502),
    @(13, 44358, 1.218113321277878),
    @(14, 44359, 0.9926902173913044),
    @(15, 44360, 1.003660130718954),
    @(16, 44361, 1.369221207508879),
    @(17, 44362, 1.372681890075129),
    @(18, 44363, 1.44971848110737),
    @(19, 44364, 1.330242782152231),
    @(20, 44365, 1.387725015518312),
    @(21, 44366, 1.414229242979243),
    @(22, 44367, 1.336816239316239),
    @(23, 44368, 1.4093259451161),
    @(24, 44369, 1.468861910471623),
    @(25, 44370, 1.40110670970266),
    @(26, 44371, 1.304626111111111),
    @(27, 44372, 1.189906859522244),
    @(28, 44373, 1.282023026315789),
    @(29, 44374, 1.371332199546485),
    @(30, 44375, 1.405848429951691),
    @(31, 44376, 1.327265715611227),
    @(32, 44377, 1.147075308641975),
    @(33, 44378, 1.166889812249923),
    @(34, 44379, 0.9453638497652582),
    @(35, 44380, 0.8778387533875339),
    @(36, 44381, 0.9814387464387465),
    @(37, 44382, 0.945027885027885),
    @(38, 44383, 0.9528630164460489),
    @(39, 44384, 0.7768076599326599),
    @(40, 44385, 1.000829365079365),
    @(41, 44386, 0.9878218482905983),
    @(42, 44387, 1.148390022675737),
    @(43, 44388, 0.8674074074074074),
    @(44, 44389, 0.9517135207496653),
    @(45, 44390, 0.724330459770115),
    @(46, 44391, 0.9171821705426356),
    @(47, 44392, 0.8364704491725768),
    @(48, 44393, 0.9037944983818771),
    @(49, 44394, 0.6472619047619048),
    @(50, 44395, 0.9096743295019157),
    @(51, 44396, 0.7815502244668912),
    @(52, 44397, 0.7386097820308346),
    @(53, 44398, 0.8869485094850948),
    @(54, 44399, 0.8601081871345029),
    @(55, 44400, 0.9018279132791328),
    @(56, 44401, 0.8836210317460317),
    @(57, 44402, 1.10858024691358),
    @(58, 44403, 0.9479040404040404),
    @(59, 44404, 0.9301051051051051),
    @(60, 44405, 1.042901678657074),
    @(61, 44406, 0.9314169215086646),
    @(62, 44407, 0.9621720430107527)
)

foreach ($row in $data) {
    $r = $row[0]
    $a = $r - 2
    $b = $row[1]
    $c = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
}
